$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.034580809857016
$ws.Cells.Item(2, 4).Value = 1.037780492126913
$ws.Cells.Item(2, 5).Value = 1.03822575358075
$ws.Cells.Item(2, 6).Value = 1.044587770365218
$ws.Cells.Item(2, 9).Value = 1.039755918629503
$ws.Cells.Item(2, 10).Value = 1.039699273166877
$ws.Cells.Item(2, 11).Value = 1.040570399656616
$ws.Cells.Item(2, 12).Value = 1.041014391730533
$ws.Cells.Item(2, 13).Value = 1.047358399524011
$ws.Cells.Item(2, 14).Value = 1.041175765389948
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.035533240766308
$ws.Cells.Item(3, 4).Value = 1.038495034182168
$ws.Cells.Item(3, 5).Value = 1.039126364437
$ws.Cells.Item(3, 6).Value = 1.045791274615673
$ws.Cells.Item(3, 9).Value = 1.040046446879131
$ws.Cells.Item(3, 10).Value = 1.040294880276605
$ws.Cells.Item(3, 11).Value = 1.041095072102385
$ws.Cells.Item(3, 12).Value = 1.041724732698938
$ws.Cells.Item(3, 13).Value = 1.048372149208502
$ws.Cells.Item(3, 14).Value = 1.041772218330088
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.036149713575338
$ws.Cells.Item(4, 4).Value = 1.038957486580054
$ws.Cells.Item(4, 5).Value = 1.039709648816163
$ws.Cells.Item(4, 6).Value = 1.046570541138812
$ws.Cells.Item(4, 9).Value = 1.040233201468797
$ws.Cells.Item(4, 10).Value = 1.040679853774693
$ws.Cells.Item(4, 11).Value = 1.041433995033528
$ws.Cells.Item(4, 12).Value = 1.042184268249752
$ws.Cells.Item(4, 13).Value = 1.049028064688214
$ws.Cells.Item(4, 14).Value = 1.042157738534701
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.036408922904061
$ws.Cells.Item(5, 4).Value = 1.039151923610618
$ws.Cells.Item(5, 5).Value = 1.039954987385247
$ws.Cells.Item(5, 6).Value = 1.04689826936627
$ws.Cells.Item(5, 9).Value = 1.040311416735644
$ws.Cells.Item(5, 10).Value = 1.040841594585946
$ws.Cells.Item(5, 11).Value = 1.041576339961814
$ws.Cells.Item(5, 12).Value = 1.042377431814377
$ws.Cells.Item(5, 13).Value = 1.049303800186537
$ws.Cells.Item(5, 14).Value = 1.042319709036458
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.036452447867579
$ws.Cells.Item(6, 4).Value = 1.039184571709481
$ws.Cells.Item(6, 5).Value = 1.039996188164879
$ws.Cells.Item(6, 6).Value = 1.046953303694411
$ws.Cells.Item(6, 9).Value = 1.040324532035151
$ws.Cells.Item(6, 10).Value = 1.040868745580937
$ws.Cells.Item(6, 11).Value = 1.041600232176957
$ws.Cells.Item(6, 12).Value = 1.042409863331947
$ws.Cells.Item(6, 13).Value = 1.04935009672941
$ws.Cells.Item(6, 14).Value = 1.042346898588976
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.036153176970108
$ws.Cells.Item(7, 4).Value = 1.038960084571716
$ws.Cells.Item(7, 5).Value = 1.039712926548776
$ws.Cells.Item(7, 6).Value = 1.046574919769421
$ws.Cells.Item(7, 9).Value = 1.040234247749914
$ws.Cells.Item(7, 10).Value = 1.040682015364944
$ws.Cells.Item(7, 11).Value = 1.041435897597446
$ws.Cells.Item(7, 12).Value = 1.042186849410489
$ws.Cells.Item(7, 13).Value = 1.049031749123846
$ws.Cells.Item(7, 14).Value = 1.042159903194658
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.03490264949362
$ws.Cells.Item(8, 4).Value = 1.038021954550219
$ws.Cells.Item(8, 5).Value = 1.038530009229371
$ws.Cells.Item(8, 6).Value = 1.044994393068555
$ws.Cells.Item(8, 9).Value = 1.039854359871663
$ws.Cells.Item(8, 10).Value = 1.039900648890347
$ws.Cells.Item(8, 11).Value = 1.040747833932667
$ws.Cells.Item(8, 12).Value = 1.041254475643119
$ws.Cells.Item(8, 13).Value = 1.047701011111569
$ws.Cells.Item(8, 14).Value = 1.041377427090043
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.03270050675891
$ws.Cells.Item(9, 4).Value = 1.036369624428763
$ws.Cells.Item(9, 5).Value = 1.036449644362436
$ws.Cells.Item(9, 6).Value = 1.042213259154611
$ws.Cells.Item(9, 9).Value = 1.03917549166535
$ws.Cells.Item(9, 10).Value = 1.038520555310849
$ws.Cells.Item(9, 11).Value = 1.039530996915325
$ws.Cells.Item(9, 12).Value = 1.039610756039763
$ws.Cells.Item(9, 13).Value = 1.045355700812302
$ws.Cells.Item(9, 14).Value = 1.039995373619363
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.03123340810243
$ws.Cells.Item(10, 4).Value = 1.035268645017429
$ws.Cells.Item(10, 5).Value = 1.0350655283515
$ws.Cells.Item(10, 6).Value = 1.040361802673109
$ws.Cells.Item(10, 9).Value = 1.038716571941578
$ws.Cells.Item(10, 10).Value = 1.03759835256572
$ws.Cells.Item(10, 11).Value = 1.038716858673666
$ws.Cells.Item(10, 12).Value = 1.038514469008764
$ws.Cells.Item(10, 13).Value = 1.043791889086679
$ws.Cells.Item(10, 14).Value = 1.039071861240558
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.030598378873195
$ws.Cells.Item(11, 4).Value = 1.034792055696787
$ws.Cells.Item(11, 5).Value = 1.034466862739046
$ws.Cells.Item(11, 6).Value = 1.039560719113775
$ws.Cells.Item(11, 9).Value = 1.038516354028157
$ws.Cells.Item(11, 10).Value = 1.037198526496314
$ws.Cells.Item(11, 11).Value = 1.038363643086758
$ws.Cells.Item(11, 12).Value = 1.038039659210206
$ws.Cells.Item(11, 13).Value = 1.043114672489831
$ws.Cells.Item(11, 14).Value = 1.038671467372273
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.030362536183254
$ws.Cells.Item(12, 4).Value = 1.034615051360719
$ws.Cells.Item(12, 5).Value = 1.034244592378856
$ws.Cells.Item(12, 6).Value = 1.039263251849873
$ws.Cells.Item(12, 9).Value = 1.038441758641662
$ws.Cells.Item(12, 10).Value = 1.037049937627065
$ws.Cells.Item(12, 11).Value = 1.03823234012521
$ws.Cells.Item(12, 12).Value = 1.037863277643067
$ws.Cells.Item(12, 13).Value = 1.042863112214053
$ws.Cells.Item(12, 14).Value = 1.038522667489786
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.030413123632892
$ws.Cells.Item(13, 4).Value = 1.034653018375842
$ws.Cells.Item(13, 5).Value = 1.034292265557625
$ws.Cells.Item(13, 6).Value = 1.039327055488398
$ws.Cells.Item(13, 9).Value = 1.038457769807896
$ws.Cells.Item(13, 10).Value = 1.037081813861717
$ws.Cells.Item(13, 11).Value = 1.038260509708337
$ws.Cells.Item(13, 12).Value = 1.037901112806978
$ws.Cells.Item(13, 13).Value = 1.042917073292518
$ws.Cells.Item(13, 14).Value = 1.038554588992348
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.030578883318387
$ws.Cells.Item(14, 4).Value = 1.034777424006406
$ws.Cells.Item(14, 5).Value = 1.034448487743462
$ws.Cells.Item(14, 6).Value = 1.039536128527928
$ws.Cells.Item(14, 9).Value = 1.038510192549902
$ws.Cells.Item(14, 10).Value = 1.037186245635337
$ws.Cells.Item(14, 11).Value = 1.038352791640848
$ws.Cells.Item(14, 12).Value = 1.038025079779705
$ws.Cells.Item(14, 13).Value = 1.043093878684102
$ws.Cells.Item(14, 14).Value = 1.038659169071064
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.030681018035392
$ws.Cells.Item(15, 4).Value = 1.034854077381897
$ws.Cells.Item(15, 5).Value = 1.034544754752777
$ws.Cells.Item(15, 6).Value = 1.039664957343834
$ws.Cells.Item(15, 9).Value = 1.038542462051858
$ws.Cells.Item(15, 10).Value = 1.037250579482964
$ws.Cells.Item(15, 11).Value = 1.038409635954735
$ws.Cells.Item(15, 12).Value = 1.038101457808992
$ws.Cells.Item(15, 13).Value = 1.043202812732539
$ws.Cells.Item(15, 14).Value = 1.038723594280134
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.031275557856701
$ws.Cells.Item(16, 4).Value = 1.035300277743424
$ws.Cells.Item(16, 5).Value = 1.035105273901415
$ws.Cells.Item(16, 6).Value = 1.040414980727971
$ws.Cells.Item(16, 9).Value = 1.038729828104549
$ws.Cells.Item(16, 10).Value = 1.037624877073075
$ws.Cells.Item(16, 11).Value = 1.038740285964601
$ws.Cells.Item(16, 12).Value = 1.038545978290656
$ws.Cells.Item(16, 13).Value = 1.043836832082869
$ws.Cells.Item(16, 14).Value = 1.039098423415756
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.031648559855274
$ws.Cells.Item(17, 4).Value = 1.035580206046332
$ws.Cells.Item(17, 5).Value = 1.03545705147739
$ws.Cells.Item(17, 6).Value = 1.040885613214028
$ws.Cells.Item(17, 9).Value = 1.038846955731132
$ws.Cells.Item(17, 10).Value = 1.037859528811651
$ws.Cells.Item(17, 11).Value = 1.038947510006953
$ws.Cells.Item(17, 12).Value = 1.038824785073573
$ws.Cells.Item(17, 13).Value = 1.044234515029661
$ws.Cells.Item(17, 14).Value = 1.039333408386715
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.031866148186898
$ws.Cells.Item(18, 4).Value = 1.035743497138399
$ws.Cells.Item(18, 5).Value = 1.035662301575762
$ws.Cells.Item(18, 6).Value = 1.04116018401774
$ws.Cells.Item(18, 9).Value = 1.038915129296051
$ws.Cells.Item(18, 10).Value = 1.03799634833442
$ws.Cells.Item(18, 11).Value = 1.039068313838249
$ws.Cells.Item(18, 12).Value = 1.038987397678998
$ws.Cells.Item(18, 13).Value = 1.044466469483345
$ws.Cells.Item(18, 14).Value = 1.0394704222089
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.031940344006525
$ws.Cells.Item(19, 4).Value = 1.035799177471119
$ws.Cells.Item(19, 5).Value = 1.035732297422141
$ws.Cells.Item(19, 6).Value = 1.041253815567539
$ws.Cells.Item(19, 9).Value = 1.038938350119177
$ws.Cells.Item(19, 10).Value = 1.038042991955655
$ws.Cells.Item(19, 11).Value = 1.03910949353861
$ws.Cells.Item(19, 12).Value = 1.039042842549416
$ws.Cells.Item(19, 13).Value = 1.044545558737821
$ws.Cells.Item(19, 14).Value = 1.039517132069427
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.031608537933727
$ws.Cells.Item(20, 4).Value = 1.035550170970895
$ws.Cells.Item(20, 5).Value = 1.035419302428224
$ws.Cells.Item(20, 6).Value = 1.040835112721428
$ws.Cells.Item(20, 9).Value = 1.038834404041507
$ws.Cells.Item(20, 10).Value = 1.037834357952458
$ws.Cells.Item(20, 11).Value = 1.03892528369793
$ws.Cells.Item(20, 12).Value = 1.038794872844748
$ws.Cells.Item(20, 13).Value = 1.044191848169608
$ws.Cells.Item(20, 14).Value = 1.039308201782016
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.030530070242036
$ws.Cells.Item(21, 4).Value = 1.034740789030701
$ws.Cells.Item(21, 5).Value = 1.034402481418074
$ws.Cells.Item(21, 6).Value = 1.039474559226848
$ws.Cells.Item(21, 9).Value = 1.038494761583209
$ws.Cells.Item(21, 10).Value = 1.037155495166522
$ws.Cells.Item(21, 11).Value = 1.038325619750516
$ws.Cells.Item(21, 12).Value = 1.037988575035534
$ws.Cells.Item(21, 13).Value = 1.043041814234866
$ws.Cells.Item(21, 14).Value = 1.038628374933056
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029852199978052
$ws.Cells.Item(22, 4).Value = 1.034232027082796
$ws.Cells.Item(22, 5).Value = 1.033763748670526
$ws.Cells.Item(22, 6).Value = 1.038619650461576
$ws.Cells.Item(22, 9).Value = 1.038279910100086
$ws.Cells.Item(22, 10).Value = 1.036728229606928
$ws.Cells.Item(22, 11).Value = 1.037947991722529
$ws.Cells.Item(22, 12).Value = 1.037481530890459
$ws.Cells.Item(22, 13).Value = 1.042318673289438
$ws.Cells.Item(22, 14).Value = 1.03820050260736
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.030211532294999
$ws.Cells.Item(23, 4).Value = 1.034501718954412
$ws.Cells.Item(23, 5).Value = 1.034102297497905
$ws.Cells.Item(23, 6).Value = 1.039072804223052
$ws.Cells.Item(23, 9).Value = 1.038393930542626
$ws.Cells.Item(23, 10).Value = 1.036954772513146
$ws.Cells.Item(23, 11).Value = 1.038148235781069
$ws.Cells.Item(23, 12).Value = 1.037750333205809
$ws.Cells.Item(23, 13).Value = 1.042702030686874
$ws.Cells.Item(23, 14).Value = 1.03842736723049
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.031626622047517
$ws.Cells.Item(24, 4).Value = 1.035563742486509
$ws.Cells.Item(24, 5).Value = 1.035436359400323
$ws.Cells.Item(24, 6).Value = 1.040857931538187
$ws.Cells.Item(24, 9).Value = 1.03884007605776
$ws.Cells.Item(24, 10).Value = 1.037845731731224
$ws.Cells.Item(24, 11).Value = 1.038935327015558
$ws.Cells.Item(24, 12).Value = 1.038808388926364
$ws.Cells.Item(24, 13).Value = 1.044211127509356
$ws.Cells.Item(24, 14).Value = 1.039319591712851
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.033269639564824
$ws.Cells.Item(25, 4).Value = 1.036796694074673
$ws.Cells.Item(25, 5).Value = 1.03698697949145
$ws.Cells.Item(25, 6).Value = 1.042931782919843
$ws.Cells.Item(25, 9).Value = 1.039352114291282
$ws.Cells.Item(25, 10).Value = 1.038877721931421
$ws.Cells.Item(25, 11).Value = 1.039846094323061
$ws.Cells.Item(25, 12).Value = 1.040035783102453
$ws.Cells.Item(25, 13).Value = 1.045962066432001
$ws.Cells.Item(25, 14).Value = 1.040353047457504

Write-Host "applied 380 kV case updates"
